$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 4: Inscritos 17 -> 18
$ws.Range("E4").Value = 18

# Row 6: Inscritos 49 -> 51
$ws.Range("E6").Value = 51

# Row 11: Inscritos 12 -> 13, Pagos 8 -> 9, Inscricoes homologadas 8 -> 9
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = 9
$ws.Range("H11").Value = 9

# Row 15: Inscritos 84 -> 85, Pagos 40 -> 41, Inscricoes homologadas 40 -> 41
$ws.Range("E15").Value = 85
$ws.Range("F15").Value = 41
$ws.Range("H15").Value = 41
